# Updated cryptos list values (Price column D, Volume(1h) column E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with the default (unstyled) look used by every data cell in column D,
# used to restore the style pointer after forcing a text NumberFormat so we don't leave
# an explicit style index behind on edited cells.
$styleDonor = $ws.Cells.Item(6, 4)

$d2 = $ws.Cells.Item(2, 4)
$d2.NumberFormat = "@"
$d2.Value = "26.055.00"
$d2.Style = $styleDonor.Style
$ws.Cells.Item(2, 5).Value = "  -0.20%  "

$d3 = $ws.Cells.Item(3, 4)
$d3.NumberFormat = "@"
$d3.Value = "1.650.56"
$d3.Style = $styleDonor.Style
$ws.Cells.Item(3, 5).Value = "  +0.04%  "

$d4 = $ws.Cells.Item(4, 4)
$d4.NumberFormat = "@"
$d4.Value = "1.003"
$d4.Style = $styleDonor.Style
$ws.Cells.Item(4, 5).Value = "  -0.23%  "

$d5 = $ws.Cells.Item(5, 4)
$d5.NumberFormat = "@"
$d5.Value = "218.32"
$d5.Style = $styleDonor.Style
$ws.Cells.Item(5, 5).Value = "  +0.41%  "

$ws.Cells.Item(6, 5).Value = "  +0.36%  "

$ws.Cells.Item(7, 5).Value = "  -0.24%  "

$d8 = $ws.Cells.Item(8, 4)
$d8.NumberFormat = "@"
$d8.Value = "0.2639"
$d8.Style = $styleDonor.Style
$ws.Cells.Item(8, 5).Value = "  +0.84%  "

$d9 = $ws.Cells.Item(9, 4)
$d9.NumberFormat = "@"
$d9.Value = "0.06325"
$d9.Style = $styleDonor.Style
$ws.Cells.Item(9, 5).Value = "  +0.75%  "

$ws.Cells.Item(10, 5).Value = "  -0.29%  "

$d11 = $ws.Cells.Item(11, 4)
$d11.NumberFormat = "@"
$d11.Value = "0.07659"
$d11.Style = $styleDonor.Style
$ws.Cells.Item(11, 5).Value = "  -1.72%  "

$d12 = $ws.Cells.Item(12, 4)
$d12.NumberFormat = "@"
$d12.Value = "4.586"
$d12.Style = $styleDonor.Style
$ws.Cells.Item(12, 5).Value = "  +2.75%  "

$d13 = $ws.Cells.Item(13, 4)
$d13.NumberFormat = "@"
$d13.Value = "1.658.98"
$d13.Style = $styleDonor.Style
$ws.Cells.Item(13, 5).Value = "  -0.58%  "

$d14 = $ws.Cells.Item(14, 4)
$d14.NumberFormat = "@"
$d14.Value = "1.877.32"
$d14.Style = $styleDonor.Style
$ws.Cells.Item(14, 5).Value = "  +0.03%  "

$d15 = $ws.Cells.Item(15, 4)
$d15.NumberFormat = "@"
$d15.Value = "0.5596"
$d15.Style = $styleDonor.Style
$ws.Cells.Item(15, 5).Value = "  +0.95%  "

$d16 = $ws.Cells.Item(16, 4)
$d16.NumberFormat = "@"
$d16.Value = "0.0₅8134"
$d16.Style = $styleDonor.Style
$ws.Cells.Item(16, 5).Value = "  +1.97%  "

$d17 = $ws.Cells.Item(17, 4)
$d17.NumberFormat = "@"
$d17.Value = "65.22"
$d17.Style = $styleDonor.Style
$ws.Cells.Item(17, 5).Value = "  +0.82%  "

$d18 = $ws.Cells.Item(18, 4)
$d18.NumberFormat = "@"
$d18.Value = "26.043.98"
$d18.Style = $styleDonor.Style
$ws.Cells.Item(18, 5).Value = "  -0.21%  "

$ws.Cells.Item(19, 5).Value = "  -0.27%  "

$d20 = $ws.Cells.Item(20, 4)
$d20.NumberFormat = "@"
$d20.Value = "4.619"
$d20.Style = $styleDonor.Style
$ws.Cells.Item(20, 5).Value = "  -0.12%  "

$d21 = $ws.Cells.Item(21, 4)
$d21.NumberFormat = "@"
$d21.Value = "10.50"
$d21.Style = $styleDonor.Style
$ws.Cells.Item(21, 5).Value = "  +4.30%  "

$d22 = $ws.Cells.Item(22, 4)
$d22.NumberFormat = "@"
$d22.Value = "191.56"
$d22.Style = $styleDonor.Style
$ws.Cells.Item(22, 5).Value = "  -1.22%  "

$d23 = $ws.Cells.Item(23, 4)
$d23.NumberFormat = "@"
$d23.Value = "5.901"
$d23.Style = $styleDonor.Style
$ws.Cells.Item(23, 5).Value = "  -0.71%  "

$ws.Cells.Item(24, 5).Value = "  -0.22%  "

$d25 = $ws.Cells.Item(25, 4)
$d25.NumberFormat = "@"
$d25.Value = "144.00"
$d25.Style = $styleDonor.Style
$ws.Cells.Item(25, 5).Value = "  -1.88%  "

$d26 = $ws.Cells.Item(26, 4)
$d26.NumberFormat = "@"
$d26.Value = "0.1185"
$d26.Style = $styleDonor.Style
$ws.Cells.Item(26, 5).Value = "  -1.65%  "

$d27 = $ws.Cells.Item(27, 4)
$d27.NumberFormat = "@"
$d27.Value = "7.193"
$d27.Style = $styleDonor.Style
$ws.Cells.Item(27, 5).Value = "  +0.45%  "

$d28 = $ws.Cells.Item(28, 4)
$d28.NumberFormat = "@"
$d28.Value = "15.88"
$d28.Style = $styleDonor.Style
$ws.Cells.Item(28, 5).Value = "  -0.27%  "

$d29 = $ws.Cells.Item(29, 4)
$d29.NumberFormat = "@"
$d29.Value = "1.514"
$d29.Style = $styleDonor.Style
$ws.Cells.Item(29, 5).Value = "  +2.58%  "

$d30 = $ws.Cells.Item(30, 4)
$d30.NumberFormat = "@"
$d30.Value = "0.05428"
$d30.Style = $styleDonor.Style
$ws.Cells.Item(30, 5).Value = "  -3.43%  "

$ws.Cells.Item(32, 5).Value = "  -1.09%  "

$ws.Cells.Item(33, 5).Value = "  -0.93%  "

$d34 = $ws.Cells.Item(34, 4)
$d34.NumberFormat = "@"
$d34.Value = "1.553"
$d34.Style = $styleDonor.Style
$ws.Cells.Item(34, 5).Value = "  -2.70%  "

$ws.Cells.Item(35, 5).Value = "  +0.80%  "

$d36 = $ws.Cells.Item(36, 4)
$d36.NumberFormat = "@"
$d36.Value = "2.781"
$d36.Style = $styleDonor.Style

$d37 = $ws.Cells.Item(37, 4)
$d37.NumberFormat = "@"
$d37.Value = "0.9444"
$d37.Style = $styleDonor.Style
$ws.Cells.Item(37, 5).Value = "  -0.27%  "

$d38 = $ws.Cells.Item(38, 4)
$d38.NumberFormat = "@"
$d38.Value = "0.5631"
$d38.Style = $styleDonor.Style
$ws.Cells.Item(38, 5).Value = "  -0.26%  "

$d39 = $ws.Cells.Item(39, 4)
$d39.NumberFormat = "@"
$d39.Value = "0.01581"
$d39.Style = $styleDonor.Style
$ws.Cells.Item(39, 5).Value = "  +0.24%  "

$ws.Cells.Item(40, 5).Value = "  -1.62%  "

$ws.Cells.Item(41, 5).Value = "  -0.17%  "

$d42 = $ws.Cells.Item(42, 4)
$d42.NumberFormat = "@"
$d42.Value = "1.025.47"
$d42.Style = $styleDonor.Style
$ws.Cells.Item(42, 5).Value = "  -3.31%  "

$d43 = $ws.Cells.Item(43, 4)
$d43.NumberFormat = "@"
$d43.Value = "0.8243"
$d43.Style = $styleDonor.Style
$ws.Cells.Item(43, 5).Value = "  -1.60%  "

$d44 = $ws.Cells.Item(44, 4)
$d44.NumberFormat = "@"
$d44.Value = "100.72"
$d44.Style = $styleDonor.Style
$ws.Cells.Item(44, 5).Value = "  -2.17%  "

$d45 = $ws.Cells.Item(45, 4)
$d45.NumberFormat = "@"
$d45.Value = "1.786.35"
$d45.Style = $styleDonor.Style
$ws.Cells.Item(45, 5).Value = "  -0.07%  "

$ws.Cells.Item(46, 5).Value = "  +5.19%  "

$d47 = $ws.Cells.Item(47, 4)
$d47.NumberFormat = "@"
$d47.Value = "57.25"
$d47.Style = $styleDonor.Style
$ws.Cells.Item(47, 5).Value = "  +0.38%  "

$d48 = $ws.Cells.Item(48, 4)
$d48.NumberFormat = "@"
$d48.Value = "0.9981"
$d48.Style = $styleDonor.Style
$ws.Cells.Item(48, 5).Value = "  -0.67%  "

$d49 = $ws.Cells.Item(49, 4)
$d49.NumberFormat = "@"
$d49.Value = "0.4332"
$d49.Style = $styleDonor.Style
$ws.Cells.Item(49, 5).Value = "  -0.03%  "

$d50 = $ws.Cells.Item(50, 4)
$d50.NumberFormat = "@"
$d50.Value = "7.956"
$d50.Style = $styleDonor.Style
$ws.Cells.Item(50, 5).Value = "  +0.07%  "

$ws.Cells.Item(51, 5).Value = "  -3.32%  "
